$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 7: "add indexed joint indicies back in" / 15 / "will not be dependant on only having 4 joints"
$ws.Cells.Item(7, 1).Value = "add indexed joint indicies back in"
$ws.Cells.Item(7, 2).Value = 15
$ws.Cells.Item(7, 3).Value = "will not be dependant on only having 4 joints"

# --- Row 8: point value 20 -> 15
$ws.Cells.Item(8, 2).Value = 15

# --- Row 9: point value 25 -> 15
$ws.Cells.Item(9, 2).Value = 15

# --- Insert two new rows after row 10 (before the old "fix memory leaks" row), inheriting
# --- row 10's formatting (style s=1 / s=9 / s=1), matching the target layout.
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(11).Insert()

# New row 12: "interpolate between 2 animations"
$ws.Cells.Item(12, 1).Value = "interpolate between 2 animations"
$ws.Cells.Item(12, 2).Value = 15

# --- Row 10: "add in normal mapping" -> "add in normal mapping with specular map", point value 20 -> 15
$ws.Cells.Item(10, 1).Value = "add in normal mapping with specular map"
$ws.Cells.Item(10, 2).Value = 15

# New row 11: "be able to roate model"
$ws.Cells.Item(11, 1).Value = "be able to roate model"
$ws.Cells.Item(11, 2).Value = 15

# --- Old row 11 ("fix memory leaks and no warnings") is now row 13; point value 15 -> 10
$ws.Cells.Item(13, 2).Value = 10

# --- Old row 13 ("Total:") is now row 15; give it a SUM formula over the rubric rows
$ws.Cells.Item(15, 2).Formula = "=SUM(B7:B13)"

# --- Update the sheet's active selection to C10, matching the author's saved cursor position
$ws.Range("C10").Select()
